$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string values used in rows 29-31
$s100 = 'leant to dive deeper within a project from architecture perspective, also social context perspective'
$s101 = 'Instead of doing more research on class level and feature level of a project, we learnt to do more work based on an architectural perspective. This will give us a more concreted sense of the whole program. Also we learnt to know the importance of looking at the social context of a project before deciding to contribute. It might be frustrated if maintainers do not maintain this project any more when you contribute.'
$s102 = '13:30 - 17:00'
$s103 = 'Soobin'
$s104 = 'finish homework3'
$s105 = 'talked about interesting open issues and PRs, finished social context of our program'
$s106 = 'It is hard to define what is interesting and also some of these PRs are quite ridiculous. Some people they just want to be a contributor by doing minimal, like fixing typo in commands. So funny and things also happen in reality.'
$s107 = 'Nervous'
$s108 = '14:00 - 17:00 & 19:30 - 21:00'
$s109 = 'finished architecture recovering'
$s110 = 'What we did in class by dragging and grouping Pacman classes in UML diagram is kinda a bottom up comprehension. But when it comes to a 100K LOC program, it is impossible to use that strategy since there are SOOOO many classes and interfaces. So we used a top down way to comprehend all the meaningful features first then top down confirm our hypothesis by taking a look at these folders’ name, files implementations and it works fine'
$s111 = 'Good, relief'

# Clone formatting/structure of row 28 (the last filled diary row) into rows 29-31
$templateRow = $ws.Range("A28:G28")
$templateRow.Copy($ws.Range("A29:G29"))
$templateRow.Copy($ws.Range("A30:G30"))
$templateRow.Copy($ws.Range("A31:G31"))

# Row 29: 2/20/2020 - Follow the lecture with professor
$ws.Range("A29").Value = 43881
$ws.Range("B29").Value = "17:00-19:00 in class"
$ws.Range("C29").Value = "N.A."
$ws.Range("D29").Value = "Follow the lecture with professor"
$ws.Range("E29").Value = $s100
$ws.Range("F29").Value = $s101
$ws.Range("G29").Value = "Average"
# D29 uses the distinct format applied to other "Follow the lecture" rows (e.g. D19)
$ws.Range("D19").Copy()
$ws.Range("D29").PasteSpecial(-4122)

# Row 30: 2/25/2020
$ws.Range("A30").Value = 43886
$ws.Range("B30").Value = $s102
$ws.Range("C30").Value = $s103
$ws.Range("D30").Value = $s104
$ws.Range("E30").Value = $s105
$ws.Range("F30").Value = $s106
$ws.Range("G30").Value = $s107

# Row 31: 2/26/2020
$ws.Range("A31").Value = 43887
$ws.Range("B31").Value = $s108
$ws.Range("C31").Value = $s103
$ws.Range("D31").Value = $s104
$ws.Range("E31").Value = $s109
$ws.Range("F31").Value = $s110
$ws.Range("G31").Value = $s111

$excel.CutCopyMode = $false

Write-Output "Updated rows 29-31 with new diary entries"
